# Registration.xlsx edit: fixed "element not found" test data
#
# 1) Update the sample-record row on "Registration-valid" to the new
#    (intentionally malformed / refreshed) test values.
# 2) Move the cell selection on "Order" from D2 to A2, and leave that
#    sheet no longer the active tab.
# 3) Make "Registration-valid" the active/selected sheet (it keeps its
#    existing D2 selection).

$wb = $excel.ActiveWorkbook

# --- 1. New sample data on Registration-valid -----------------------------
$wsReg = $wb.Worksheets.Item("Registration-valid")
$wsReg.Range("A2").Value = "Fardinddppp"
$wsReg.Range("B2").Value = "Akjnjkhosandd"
$wsReg.Range("C2").Value = "fardinahhosan@gmail.comd"
$wsReg.Range("D2").Value = "uyeguyegcuecvwcellgvdd"

# --- 2. Move the Order sheet's selection to A2 -----------------------------
$wsOrder = $wb.Worksheets.Item("Order")
$wsOrder.Activate()
[void]$wsOrder.Range("A2").Select()

# --- 3. Re-select Registration-valid so it becomes the active tab ---------
$wsReg.Activate()
[void]$wsReg.Range("D2").Select()
